$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.2
$ws.Range("M2").Value = 4.5
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 1.11
$ws.Range("K3").Value = 6.5
$ws.Range("L3").Value = 1.53
$ws.Range("M3").Value = 2.5
$ws.Range("N3").Value = 2.6
$ws.Range("O3").Value = 1.48
$ws.Range("R3").Value = 2.38
$ws.Range("S3").Value = 1.53
$ws.Range("U3").Value = 6.5
$ws.Range("W3").Value = 13
$ws.Range("Z3").Value = 6.5
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 101
$ws.Range("AF3").Value = 23
$ws.Range("AG3").Value = 19
$ws.Range("AI3").Value = 51
$ws.Range("I5").Value = 3.75
$ws.Range("L5").Value = 1.36
$ws.Range("M5").Value = 3.2
$ws.Range("Z5").Value = 8.5
$ws.Range("AI5").Value = 34
$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 4
$ws.Range("R6").Value = 2.2
$ws.Range("S6").Value = 1.62
$ws.Range("U6").Value = 8.5
$ws.Range("V6").Value = 10
$ws.Range("W6").Value = 19
$ws.Range("N8").Value = 1.62
$ws.Range("O8").Value = 2.25
$ws.Range("J9").Value = 1.1
$ws.Range("K9").Value = 7
$ws.Range("T9").Value = 6.5
$ws.Range("X9").Value = 23
$ws.Range("AD9").Value = 451
$ws.Range("AI9").Value = 29
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.8
$ws.Range("P11").Value = 1.67
$ws.Range("Q11").Value = 2.1
$ws.Range("T11").Value = 5.5
$ws.Range("U11").Value = 8.5
$ws.Range("W11").Value = 21
$ws.Range("X11").Value = 23
$ws.Range("AF11").Value = 17
$ws.Range("G12").Value = 1.73
$ws.Range("H12").Value = 3.7
$ws.Range("J12").Value = 1.07
$ws.Range("K12").Value = 9
$ws.Range("AE12").Value = 11
$ws.Range("G13").Value = 2.5
$ws.Range("I13").Value = 3.25
$ws.Range("R13").Value = 2.38
$ws.Range("S13").Value = 1.53
$ws.Range("T13").Value = 5.5
$ws.Range("U13").Value = 10
$ws.Range("Z13").Value = 5
$ws.Range("AA13").Value = 5.5
$ws.Range("AF13").Value = 15
$ws.Range("G20").Value = 1.65
$ws.Range("H20").Value = 3.4
$ws.Range("J20").Value = 1.1
$ws.Range("K20").Value = 7
$ws.Range("L20").Value = 1.44
$ws.Range("M20").Value = 2.63
$ws.Range("N20").Value = 2.5
$ws.Range("O20").Value = 1.5
$ws.Range("P20").Value = 1.53
$ws.Range("Q20").Value = 2.38
$ws.Range("R20").Value = 2.38
$ws.Range("S20").Value = 1.53
$ws.Range("T20").Value = 5
$ws.Range("W20").Value = 12
$ws.Range("Z20").Value = 6.5
$ws.Range("AB20").Value = 23
$ws.Range("AE20").Value = 11
$ws.Range("AJ20").Value = 67
$ws.Range("G21").Value = 1.2
$ws.Range("H21").Value = 6
$ws.Range("I21").Value = 15
$ws.Range("V21").Value = 10
$ws.Range("Y21").Value = 34
$ws.Range("AB21").Value = 29
$ws.Range("G24").Value = 2.1
$ws.Range("I24").Value = 3.25
$ws.Range("N24").Value = 2.03
$ws.Range("O24").Value = 1.78
$ws.Range("N41").Value = 1.88
$ws.Range("O41").Value = 1.98
$ws.Range("W41").Value = 26
$ws.Range("Z41").Value = 11
$ws.Range("G42").Value = 3.4
$ws.Range("I42").Value = 2.2
$ws.Range("L42").Value = 1.3
$ws.Range("M42").Value = 3.4
$ws.Range("T42").Value = 10
$ws.Range("X42").Value = 29
$ws.Range("Z42").Value = 9
$ws.Range("AF42").Value = 10
$ws.Range("AG42").Value = 9
$ws.Range("P43").Value = 1.25
$ws.Range("Q43").Value = 3.75
$ws.Range("R43").Value = 1.5
$ws.Range("S43").Value = 2.5
$ws.Range("AD43").Value = 101
$ws.Range("G45").Value = 2.6
$ws.Range("H45").Value = 3
$ws.Range("I45").Value = 2.88
$ws.Range("J45").Value = 1.11
$ws.Range("K45").Value = 6.5
$ws.Range("N45").Value = 2.6
$ws.Range("O45").Value = 1.48
$ws.Range("T45").Value = 6.5
$ws.Range("U45").Value = 11
$ws.Range("V45").Value = 11
$ws.Range("W45").Value = 26
$ws.Range("X45").Value = 26
$ws.Range("AE45").Value = 7
$ws.Range("AF45").Value = 12
$ws.Range("AG45").Value = 12
$ws.Range("AH45").Value = 29
$ws.Range("AI45").Value = 29
$ws.Range("AJ45").Value = 41
